# Update the "想去人数" (interested-count) figures in the two sheets that
# list the full event data ("展览" and "全部类型") to reflect newly
# generated output data (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F6").Value = 11100
    $ws.Range("F7").Value = 561
    $ws.Range("F19").Value = 1182
}
